# Update "想去人数" (interest count) values in the F column for the
# sheets that contain conference data: "展览" and "全部类型".
# Row -> (old value -> new value):
#   F2  : 72   -> 73
#   F3  : 1061 -> 1062
#   F5  : 3053 -> 3052
#   F7  : 2196 -> 2202
#   F9  : 108  -> 110
#   F10 : 1033 -> 1039
#   F13 : 253  -> 256
#   F16 : 43   -> 44

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    2  = 73
    3  = 1062
    5  = 3052
    7  = 2202
    9  = 110
    10 = 1039
    13 = 256
    16 = 44
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
